$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BANBK")

# Row 7
$ws.Range("F7").Value = 236.8
$ws.Range("G7").Value = 241
$ws.Range("H7").Value = 232.95
$ws.Range("I7").Value = 238.25
$ws.Range("J7").Value = 239.3

# Row 9
$ws.Range("G9").Value = 238.3
$ws.Range("H9").Value = 235
$ws.Range("I9").Value = 236.05

# Row 10
$ws.Range("G10").Value = 240.2
$ws.Range("H10").Value = 235.05
$ws.Range("I10").Value = 239.95

# Row 11
$ws.Range("G11").Value = 241
$ws.Range("H11").Value = 239.65
$ws.Range("I11").Value = 239.95

# Row 12
$ws.Range("G12").Value = 240.75
$ws.Range("H12").Value = 237.75
$ws.Range("I12").Value = 238

# Row 13
$ws.Range("G13").Value = 238.3
$ws.Range("H13").Value = 236.55
$ws.Range("I13").Value = 236.7

# Row 14
$ws.Range("G14").Value = 237.9
$ws.Range("H14").Value = 236.65
$ws.Range("I14").Value = 237.25

# Row 15
$ws.Range("G15").Value = 237.9
$ws.Range("H15").Value = 236.55
$ws.Range("I15").Value = 237.9

# Row 16
$ws.Range("G16").Value = 239.45
$ws.Range("H16").Value = 237.5
$ws.Range("I16").Value = 237.8

# Row 17
$ws.Range("G17").Value = 238.8
$ws.Range("H17").Value = 237
$ws.Range("I17").Value = 237.2

# Row 18
$ws.Range("G18").Value = 237.6
$ws.Range("H18").Value = 235.45
$ws.Range("I18").Value = 235.7

# Row 19
$ws.Range("G19").Value = 236.5
$ws.Range("H19").Value = 232.95
$ws.Range("I19").Value = 236.25

# Row 20
$ws.Range("G20").Value = 239.2
$ws.Range("H20").Value = 235.85
$ws.Range("I20").Value = 238.5

# Row 21
$ws.Range("G21").Value = 239
$ws.Range("H21").Value = 237.15
$ws.Range("I21").Value = 238.15
